# Update odds values in row 3 of Sheet1 to reflect the latest FlashScore data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "G3"  = 1.5
    "I3"  = 6.25
    "J3"  = 2.1
    "L3"  = 6.5
    "M3"  = 1.05
    "N3"  = 11
    "O3"  = 1.29
    "P3"  = 3.5
    "Q3"  = 1.98
    "R3"  = 1.88
    "U3"  = 2
    "V3"  = 1.73
    "W3"  = 6.5
    "Y3"  = 8.5
    "Z3"  = 10
    "AB3" = 29
    "AC3" = 10
    "AD3" = 7.5
    "AE3" = 19
    "AK3" = 67
    "AN3" = 3.4
    "AU3" = 9
    "AX3" = 7.5
    "AY3" = 34
    "BA3" = 126
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
